# 03-05-2024 OrangeHRM_Configuration Test Cases Completed
# Extends the "OrangeHRM_AllTestcases" sheet with a new block of
# "OrangeHRM_Configuration" test-case columns (U:AK) on the header row (6)
# and the data row (7), plus the supporting column widths, view state and
# the new hyperlink on AD7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrangeHRM_AllTestcases")

# ---------------------------------------------------------------------
# 1. New header row (row 6) values, columns U:AK -- same fill style (s=17)
#    as the rest of row 6, picked up automatically by copying format from
#    an existing header cell (A6) onto the new ones.
# ---------------------------------------------------------------------
$ws.Range("U6").Value = "UpdateOAuthName"
$ws.Range("V6").Value = "UpdateReqiuredURI"
$ws.Range("W6").Value = "HostName"
$ws.Range("X6").Value = "Port"
$ws.Range("Y6").Value = "UserAttribute"
$ws.Range("Z6").Value = "UserSearch"
$ws.Range("AA6").Value = "UserUnique"
$ws.Range("AB6").Value = "FirstName"
$ws.Range("AC6").Value = "Middlename"
$ws.Range("AD6").Value = "WorkemailBox"
$ws.Range("AE6").Value = "EmployIdBox"
$ws.Range("AF6").Value = "SyncBox"
$ws.Range("AG6").Value = "LastName"
$ws.Range("AH6").Value = "Distinguished"
$ws.Range("AI6").Value = "passwordname"
$ws.Range("AJ6").Value = "BaseDistinguished"
$ws.Range("AK6").Value = "StatusBox"

$ws.Range("A6").Copy()
$ws.Range("U6:AK6").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. New data row (row 7) values, columns U:AK.
#    Most are plain text; a handful (X7,Y7,AA7,AE7,AF7) were typed with a
#    leading apostrophe in the original workbook (quote-prefixed numeric-
#    looking text); AD7 is a hyperlinked e-mail address; AI7 is a real
#    number.
# ---------------------------------------------------------------------
$ws.Range("U7").Value = "madhu"
$ws.Range("V7").Value = "tgsliafeh"
$ws.Range("W7").Value = "publichost"
$ws.Range("X7").Value = "'1010"
$ws.Range("Y7").Value = "'cmd"
$ws.Range("Z7").Value = "debt"
$ws.Range("AA7").Value = "'241"
$ws.Range("AB7").Value = "ram"
$ws.Range("AC7").Value = "reddy"
$ws.Range("AD7").Value = "rtb@gmail.comm"
$ws.Range("AE7").Value = "'789"
$ws.Range("AF7").Value = "'20"
$ws.Range("AG7").Value = "raj"
$ws.Range("AH7").Value = "sekhar"
$ws.Range("AI7").Value = 123456
$ws.Range("AJ7").Value = "Ramu"
$ws.Range("AK7").Value = "pass"

# Hyperlink for the new work-email cell, matching the style already used
# by the other mailto/http hyperlink cells on this sheet (F7, H7, L7, P7).
$ws.Hyperlinks.Add($ws.Range("AD7"), "mailto:rtb@gmail.comm")
$ws.Range("F7").Copy()
$ws.Range("AD7").PasteSpecial(-4122)
$ws.Range("AD7").Value = "rtb@gmail.comm"

# ---------------------------------------------------------------------
# 3. Column widths for the newly populated columns (U:AJ), matching the
#    widths the sheet's other configured columns use. ColumnWidth is
#    expressed in characters; pick the values landing closest to the
#    target stored widths.
# ---------------------------------------------------------------------
$ws.Columns.Item(21).ColumnWidth = 17.165   # U
$ws.Columns.Item(22).ColumnWidth = 18       # V
$ws.Columns.Item(23).ColumnWidth = 11.33    # W
$ws.Columns.Item(24).ColumnWidth = 12.665   # X
$ws.Columns.Item(25).ColumnWidth = 19       # Y
$ws.Columns.Item(26).ColumnWidth = 17       # Z
$ws.Columns.Item(27).ColumnWidth = 19.33    # AA
$ws.Columns.Item(28).ColumnWidth = 14.83    # AB
$ws.Columns.Item(29).ColumnWidth = 18.5     # AC
$ws.Columns.Item(30).ColumnWidth = 20       # AD
$ws.Columns.Item(31).ColumnWidth = 20       # AE
$ws.Columns.Item(34).ColumnWidth = 14.665   # AH
$ws.Columns.Item(35).ColumnWidth = 14.83    # AI
$ws.Columns.Item(36).ColumnWidth = 20.33    # AJ

# ---------------------------------------------------------------------
# 4. View-state: scroll the window over to the new columns and move the
#    active selection, matching the saved sheet view in the edit.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Application.ActiveWindow.ScrollColumn = 22
$ws.Range("Y9").Select()
